# aggiornamento fino a 13/03 - append 4 new daily rows (252-255) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (251) as a formatting template: copying a
# range of the same shape onto new rows carries over both values and cell
# styles (including the date-format style used in column A), then we
# overwrite the copied values with the real new data.
$template = $ws.Range("A251:D251")

$template.Copy($ws.Range("A252:D252"))
$template.Copy($ws.Range("A253:D253"))
$template.Copy($ws.Range("A254:D254"))
$template.Copy($ws.Range("A255:D255"))

$ws.Range("A252").Value = 44326
$ws.Range("B252").Value = 2
$ws.Range("C252").Value = 21
$ws.Range("D252").Value = 61.10690799045568

$ws.Range("A253").Value = 44327
$ws.Range("B253").Value = 1
$ws.Range("C253").Value = 22
$ws.Range("D253").Value = 64.01676075190595

$ws.Range("A254").Value = 44328
$ws.Range("B254").Value = 4
$ws.Range("C254").Value = 24
$ws.Range("D254").Value = 69.8364662748065

$ws.Range("A255").Value = 44329
$ws.Range("B255").Value = 5
$ws.Range("C255").Value = 24
$ws.Range("D255").Value = 69.8364662748065
